$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$chars = $tr.Characters(1, 8)
Write-Host "chars text:" $chars.Text
$chars.Text = "20220202"
